$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.078.09'
$ws.Range('E2').Value = '  +0.09%  '

$ws.Range('D3').Value = '1.836.05'
$ws.Range('E3').Value = '  +0.39%  '

$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '242.91'
$ws.Range('E5').Value = '  +0.56%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6270'
$ws.Range('E6').Value = '  -0.51%  '

$ws.Range('E7').Value = '  +0.07%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07587'
$ws.Range('E8').Value = '  +3.47%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2930'
$ws.Range('E9').Value = '  -0.06%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.60'
$ws.Range('E10').Value = '  -1.20%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07746'
$ws.Range('E11').Value = '  +0.92%  '

$ws.Range('D12').Value = '1.845.10'
$ws.Range('E12').Value = '  +0.90%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.971'
$ws.Range('E13').Value = '  -0.30%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6649'
$ws.Range('E14').Value = '  +0.34%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '82.91'
$ws.Range('E15').Value = '  +1.06%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000009952'

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.066'
$ws.Range('E17').Value = '  +0.18%  '

$ws.Range('D18').Value = '29.093.72'
$ws.Range('E18').Value = '  +0.15%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '227.01'
$ws.Range('E19').Value = '  +1.40%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.41'
$ws.Range('E20').Value = '  +0.12%  '

$ws.Range('E21').Value = '  +0.17%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.219'
$ws.Range('E22').Value = '  +1.29%  '

$ws.Range('E23').Value = '  +0.07%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '159.49'
$ws.Range('E24').Value = '  +0.98%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.517'
$ws.Range('E25').Value = '  +0.83%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1384'
$ws.Range('E26').Value = '  +1.34%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.94'
$ws.Range('E27').Value = '  +0.38%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.494'
$ws.Range('E28').Value = '  -0.85%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.106'
$ws.Range('E29').Value = '  +0.49%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.022'
$ws.Range('E30').Value = '  +0.10%  '

$ws.Range('E31').Value = '  -0.58%  '

$ws.Range('E32').Value = '  -0.87%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.842'
$ws.Range('E33').Value = '  +0.82%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7345'
$ws.Range('E34').Value = '  -0.60%  '

$ws.Range('E35').Value = '  -1.10%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.700'
$ws.Range('E36').Value = '  +1.77%  '

$ws.Range('D37').Value = '1.239.72'
$ws.Range('E37').Value = '  -4.24%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.763'
$ws.Range('E38').Value = '  +0.89%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01785'
$ws.Range('E39').Value = '  +0.24%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.369'
$ws.Range('E40').Value = '  +0.81%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8997'
$ws.Range('E41').Value = '  +0.56%  '

$ws.Range('E42').Value = '  +0.14%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '101.96'
$ws.Range('E43').Value = '  -0.61%  '

$ws.Range('D44').Value = '1.986.30'
$ws.Range('E44').Value = '  +0.40%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000125'
$ws.Range('E45').Value = '  +1.10%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '64.34'
$ws.Range('E46').Value = '  +0.15%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5117'
$ws.Range('E47').Value = '  -0.42%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4040'

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.878'
$ws.Range('E49').Value = '  +1.82%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05764'
$ws.Range('E50').Value = '  -1.10%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.691'
$ws.Range('E51').Value = '  +0.11%  '
